$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Reference / Payment reference values for rows 5 and 6
$ws.Range("H5").Value = "EXI00446571"
$ws.Range("I5").Value = "EXI00446571"
$ws.Range("H6").Value = "EXI00447099"
$ws.Range("I6").Value = "EXI00447099"

# Column H got a custom width in the edit (likely from autofit after entry).
# Note: the host's ColumnWidth setter snaps to its own internal grid, so an
# input of 10.33 is the closest reachable value to the target 11.140625.
$ws.Columns("H").ColumnWidth = 10.33
